$wb = $excel.ActiveWorkbook

# ============================================================
# 1) "总计" sheet: insert a new row for the 2022-Q3 summary
#    (existing 2022-Q2 / 2021-Q1 rows shift down by one)
# ============================================================
$total = $wb.Worksheets.Item("总计")

$total.Rows("2").Insert()
# Give the new A2 the same (bold + bordered + centered) style as the
# existing index column cells, then strip any leftover row-insert
# formatting from B2:D2 so they stay plain, like the originals.
$total.Range("A4").Copy()
$total.Range("A2").PasteSpecial(-4122)
$total.Range("B2:D2").ClearFormats()

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 5
$total.Range("D2").Value = 0.81

$total.Range("A3").Value = 1
$total.Range("A4").Value = 2

# ============================================================
# 2) New "2022-Q3" worksheet: clone the "2022-Q2" sheet (so the
#    header row / index-column style match exactly), position it
#    before "2022-Q2", rename it, and replace its data.
# ============================================================
$src = $wb.Worksheets.Item("2022-Q2")
$src.Copy($src)
$new = $wb.Worksheets.Item("2022-Q2 (2)")
$new.Name = "2022-Q3"

# The clone has rows 1-3 (header + 2 data rows); grow it to 6 rows
# (header + 5 data rows), copying the data-row style onto each new row.
$new.Rows("4").Insert()
$new.Range("A3").Copy()
$new.Range("A4").PasteSpecial(-4122)
$new.Range("B4:H4").ClearFormats()

$new.Rows("5").Insert()
$new.Range("A3").Copy()
$new.Range("A5").PasteSpecial(-4122)
$new.Range("B5:H5").ClearFormats()

$new.Rows("6").Insert()
$new.Range("A3").Copy()
$new.Range("A6").PasteSpecial(-4122)
$new.Range("B6:H6").ClearFormats()

function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

# Row 2
$new.Range("A2").Value = 0
Set-TextValue $new.Range("B2") "213003"
$new.Range("C2").Value = "宝盈策略增长混合"
Set-TextValue $new.Range("D2") "10.55"
Set-TextValue $new.Range("E2") "90.74"
Set-TextValue $new.Range("F2") "5.19"
Set-TextValue $new.Range("G2") "0.5475"
$new.Range("H2").Value = 7

# Row 3
$new.Range("A3").Value = 1
Set-TextValue $new.Range("B3") "010296"
$new.Range("C3").Value = "万家互联互通中国优势量化策略混合A"
Set-TextValue $new.Range("D3") "4.22"
Set-TextValue $new.Range("E3") "86.53"
Set-TextValue $new.Range("F3") "5.28"
Set-TextValue $new.Range("G3") "0.2228"
$new.Range("H3").Value = 8

# Row 4
$new.Range("A4").Value = 2
Set-TextValue $new.Range("B4") "010297"
$new.Range("C4").Value = "万家互联互通中国优势量化策略混合C"
Set-TextValue $new.Range("D4") "0.46"
Set-TextValue $new.Range("E4") "86.53"
Set-TextValue $new.Range("F4") "5.28"
Set-TextValue $new.Range("G4") "0.0243"
$new.Range("H4").Value = 8

# Row 5
$new.Range("A5").Value = 3
Set-TextValue $new.Range("B5") "740001"
$new.Range("C5").Value = "长安宏观策略混合A"
Set-TextValue $new.Range("D5") "0.30"
Set-TextValue $new.Range("E5") "92.99"
Set-TextValue $new.Range("F5") "5.29"
Set-TextValue $new.Range("G5") "0.0159"
$new.Range("H5").Value = 5

# Row 6
$new.Range("A6").Value = 4
Set-TextValue $new.Range("B6") "016579"
$new.Range("C6").Value = "长安宏观策略混合C"
Set-TextValue $new.Range("D6") "0.03"
Set-TextValue $new.Range("E6") "92.99"
Set-TextValue $new.Range("F6") "5.29"
Set-TextValue $new.Range("G6") "0.0016"
$new.Range("H6").Value = 5
